# Refresh market-board derived profit figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# columns H:N) for the leve rows whose source data changed in this scheduled run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow / Beeswax
$ws.Cells.Item(12, 8).Value = 4998  # H12
$ws.Cells.Item(12, 9).Value = 4998  # I12
$ws.Cells.Item(12, 11).Value = 4998  # K12
$ws.Cells.Item(12, 13).Value = -4828  # M12

# Row 40: Stuck in the Moment / Horn Glue
$ws.Cells.Item(40, 8).Value = 4041.5  # H40
$ws.Cells.Item(40, 9).Value = 2944.3333  # I40
$ws.Cells.Item(40, 11).Value = 2944.3333  # K40
$ws.Cells.Item(40, 13).Value = -2769.3333  # M40

# Row 95: Official Strategy Guide / Gyuki Leather Codex
$ws.Cells.Item(95, 8).Value = 13250  # H95
$ws.Cells.Item(95, 10).Value = 13250  # J95
$ws.Cells.Item(95, 12).Value = 13250  # L95
$ws.Cells.Item(95, 14).Value = -18742  # N95

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Cells.Item(132, 8).Value = 3866.125  # H132
$ws.Cells.Item(132, 9).Value = 2099.889  # I132
$ws.Cells.Item(132, 10).Value = 6137  # J132
$ws.Cells.Item(132, 11).Value = 6299.667  # K132
$ws.Cells.Item(132, 12).Value = 18411  # L132
$ws.Cells.Item(132, 13).Value = -3769.667  # M132
$ws.Cells.Item(132, 14).Value = -23471  # N132

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 3155.4443  # H137
$ws.Cells.Item(137, 9).Value = 2099.75  # I137
$ws.Cells.Item(137, 11).Value = 6299.25  # K137
$ws.Cells.Item(137, 13).Value = -3749.25  # M137

$ws = $wb.Worksheets.Item("ARM")
# Row 14: Waste Not, Want Not / Bronze Celata
$ws.Cells.Item(14, 8).Value = 825.25  # H14
$ws.Cells.Item(14, 9).Value = 867  # I14
$ws.Cells.Item(14, 10).Value = 700  # J14
$ws.Cells.Item(14, 11).Value = 867  # K14
$ws.Cells.Item(14, 12).Value = 700  # L14
$ws.Cells.Item(14, 13).Value = -692  # M14
$ws.Cells.Item(14, 14).Value = -1050  # N14

# Row 16: Greavous Losses / Bronze Sabatons
$ws.Cells.Item(16, 8).Value = 953  # H16
$ws.Cells.Item(16, 9).Value = 906  # I16
$ws.Cells.Item(16, 10).Value = 1000  # J16
$ws.Cells.Item(16, 11).Value = 906  # K16
$ws.Cells.Item(16, 12).Value = 1000  # L16
$ws.Cells.Item(16, 13).Value = -619  # M16
$ws.Cells.Item(16, 14).Value = -1574  # N16

# Row 18: Still the Best / Brass Alembic
$ws.Cells.Item(18, 14).ClearContents()  # N18
$ws.Cells.Item(18, 8).Value = 1011  # H18
$ws.Cells.Item(18, 9).Value = 1011  # I18
$ws.Cells.Item(18, 10).Value = 0  # J18
$ws.Cells.Item(18, 11).Value = 1011  # K18
$ws.Cells.Item(18, 12).Value = 0  # L18
$ws.Cells.Item(18, 13).Value = -689  # M18

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Cells.Item(45, 8).Value = 1720  # H45
$ws.Cells.Item(45, 9).Value = 1578.8889  # I45
$ws.Cells.Item(45, 10).Value = 2990  # J45
$ws.Cells.Item(45, 11).Value = 1578.8889  # K45
$ws.Cells.Item(45, 12).Value = 2990  # L45
$ws.Cells.Item(45, 13).Value = -1201.8889  # M45
$ws.Cells.Item(45, 14).Value = -3744  # N45

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Cells.Item(102, 8).Value = 1655.2727  # H102
$ws.Cells.Item(102, 9).Value = 1120.8  # I102
$ws.Cells.Item(102, 11).Value = 1120.8  # K102
$ws.Cells.Item(102, 13).Value = 501.2  # M102

# Row 135: Forgiveness for My Shins / Ruthenium Sabatons of Fending
$ws.Cells.Item(135, 8).Value = 12549858  # H135
$ws.Cells.Item(135, 10).Value = 66476.664  # J135
$ws.Cells.Item(135, 12).Value = 66476.664  # L135
$ws.Cells.Item(135, 14).Value = -76616.664  # N135

$ws = $wb.Worksheets.Item("BSM")
# Row 80: Unbreaker / Titanium Ingot
$ws.Cells.Item(80, 8).Value = 529.1177  # H80
$ws.Cells.Item(80, 9).Value = 513.625  # I80
$ws.Cells.Item(80, 10).Value = 542.8889  # J80
$ws.Cells.Item(80, 11).Value = 513.625  # K80
$ws.Cells.Item(80, 12).Value = 542.8889  # L80
$ws.Cells.Item(80, 13).Value = 484.375  # M80
$ws.Cells.Item(80, 14).Value = -2538.8889  # N80

# Row 83: Attack on Titanium (L) / Titanium Ingot
$ws.Cells.Item(83, 8).Value = 529.1177  # H83
$ws.Cells.Item(83, 9).Value = 513.625  # I83
$ws.Cells.Item(83, 10).Value = 542.8889  # J83
$ws.Cells.Item(83, 11).Value = 2568.125  # K83
$ws.Cells.Item(83, 12).Value = 2714.4445  # L83
$ws.Cells.Item(83, 13).Value = 2423.875  # M83
$ws.Cells.Item(83, 14).Value = -12698.4445  # N83

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Cells.Item(105, 8).Value = 5606.4287  # H105
$ws.Cells.Item(105, 9).Value = 1745  # I105
$ws.Cells.Item(105, 10).Value = 6250  # J105
$ws.Cells.Item(105, 11).Value = 1745  # K105
$ws.Cells.Item(105, 12).Value = 6250  # L105
$ws.Cells.Item(105, 13).Value = 2  # M105
$ws.Cells.Item(105, 14).Value = -9744  # N105

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Cells.Item(134, 8).Value = 1505.8  # H134
$ws.Cells.Item(134, 9).Value = 1505.8  # I134
$ws.Cells.Item(134, 11).Value = 4517.4  # K134
$ws.Cells.Item(134, 13).Value = -1982.4  # M134

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall / Elm Lumber
$ws.Cells.Item(22, 8).Value = 1133.3334  # H22
$ws.Cells.Item(22, 9).Value = 501  # I22
$ws.Cells.Item(22, 10).Value = 1259.8  # J22
$ws.Cells.Item(22, 11).Value = 501  # K22
$ws.Cells.Item(22, 12).Value = 1259.8  # L22
$ws.Cells.Item(22, 13).Value = -151  # M22
$ws.Cells.Item(22, 14).Value = -1959.8  # N22

# Row 96: Composition / Larch Composite Bow
$ws.Cells.Item(96, 8).Value = 18699.2  # H96
$ws.Cells.Item(96, 10).Value = 18699.2  # J96
$ws.Cells.Item(96, 12).Value = 18699.2  # L96
$ws.Cells.Item(96, 14).Value = -24191.2  # N96

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water / Boiled Egg
$ws.Cells.Item(4, 8).Value = 500250020  # H4
$ws.Cells.Item(4, 9).Value = 666666700  # I4
$ws.Cells.Item(4, 11).Value = 2000000100  # K4
$ws.Cells.Item(4, 13).Value = -1999999988  # M4

# Row 88: Don't Let It Fall Apart / Liver-cheese Sandwich
$ws.Cells.Item(88, 8).Value = 1333.3334  # H88
$ws.Cells.Item(88, 9).Value = 1333.3334  # I88
$ws.Cells.Item(88, 11).Value = 4000.0002  # K88
$ws.Cells.Item(88, 13).Value = -3572.0002  # M88

# Row 91: Better Come Back with a Sandwich (L) / Liver-cheese Sandwich
$ws.Cells.Item(91, 8).Value = 1333.3334  # H91
$ws.Cells.Item(91, 9).Value = 1333.3334  # I91
$ws.Cells.Item(91, 11).Value = 4000.0002  # K91
$ws.Cells.Item(91, 13).Value = -2518.0002  # M91

# Row 98: Sweet Kiss of Death / Rice Vinegar
$ws.Cells.Item(98, 8).Value = 554.25  # H98
$ws.Cells.Item(98, 10).Value = 574.4286  # J98
$ws.Cells.Item(98, 12).Value = 1723.2858  # L98
$ws.Cells.Item(98, 14).Value = -4719.2858  # N98

# Row 107: Slippery Service / Frantoio Oil
$ws.Cells.Item(107, 8).Value = 1597.5834  # H107
$ws.Cells.Item(107, 9).Value = 580  # I107
$ws.Cells.Item(107, 10).Value = 1865.3684  # J107
$ws.Cells.Item(107, 11).Value = 1740  # K107
$ws.Cells.Item(107, 12).Value = 5596.1052  # L107
$ws.Cells.Item(107, 13).Value = 180  # M107
$ws.Cells.Item(107, 14).Value = -9436.1052  # N107

# Row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Cells.Item(121, 8).Value = 1874.5  # H121
$ws.Cells.Item(121, 9).Value = 749.5  # I121
$ws.Cells.Item(121, 10).Value = 2999.5  # J121
$ws.Cells.Item(121, 11).Value = 2248.5  # K121
$ws.Cells.Item(121, 12).Value = 8998.5  # L121
$ws.Cells.Item(121, 13).Value = -938.5  # M121
$ws.Cells.Item(121, 14).Value = -11618.5  # N121

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Cells.Item(131, 8).Value = 1466.5  # H131
$ws.Cells.Item(131, 9).Value = 932.5  # I131
$ws.Cells.Item(131, 11).Value = 2797.5  # K131
$ws.Cells.Item(131, 13).Value = 2242.5  # M131

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Cells.Item(70, 8).Value = 7700137  # H70
$ws.Cells.Item(70, 10).Value = 4124  # J70
$ws.Cells.Item(70, 12).Value = 4124  # L70
$ws.Cells.Item(70, 14).Value = -4664  # N70

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Cells.Item(73, 8).Value = 7700137  # H73
$ws.Cells.Item(73, 10).Value = 4124  # J73
$ws.Cells.Item(73, 12).Value = 4124  # L73
$ws.Cells.Item(73, 14).Value = -5996  # N73

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Cells.Item(7, 8).Value = 4332.6665  # H7
$ws.Cells.Item(7, 9).Value = 4399.8  # I7
$ws.Cells.Item(7, 10).Value = 3997  # J7
$ws.Cells.Item(7, 11).Value = 4399.8  # K7
$ws.Cells.Item(7, 12).Value = 3997  # L7
$ws.Cells.Item(7, 13).Value = -4287.8  # M7
$ws.Cells.Item(7, 14).Value = -4221  # N7

# Row 40: Best Served Toad / Toad Leather
$ws.Cells.Item(40, 8).Value = 6589.7  # H40
$ws.Cells.Item(40, 9).Value = 4655.222  # I40
$ws.Cells.Item(40, 10).Value = 24000  # J40
$ws.Cells.Item(40, 11).Value = 4655.222  # K40
$ws.Cells.Item(40, 12).Value = 24000  # L40
$ws.Cells.Item(40, 13).Value = -4519.222  # M40
$ws.Cells.Item(40, 14).Value = -24272  # N40

# Row 46: Supply Side Logic / Boar Leather
$ws.Cells.Item(46, 8).Value = 3444.6667  # H46
$ws.Cells.Item(46, 10).Value = 3444.6667  # J46
$ws.Cells.Item(46, 12).Value = 3444.6667  # L46
$ws.Cells.Item(46, 14).Value = -3820.6667  # N46

# Row 126: Battered Books / Saiga Leather
$ws.Cells.Item(126, 8).Value = 4332.6665  # H126
$ws.Cells.Item(126, 9).Value = 4399.8  # I126
$ws.Cells.Item(126, 10).Value = 3997  # J126
$ws.Cells.Item(126, 11).Value = 13199.4  # K126
$ws.Cells.Item(126, 12).Value = 11991  # L126
$ws.Cells.Item(126, 13).Value = -10729.4  # M126
$ws.Cells.Item(126, 14).Value = -16931  # N126

$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table / Pixie Floss
$ws.Cells.Item(113, 14).ClearContents()  # N113
$ws.Cells.Item(113, 8).Value = 972.5  # H113
$ws.Cells.Item(113, 9).Value = 972.5  # I113
$ws.Cells.Item(113, 10).Value = 0  # J113
$ws.Cells.Item(113, 11).Value = 2917.5  # K113
$ws.Cells.Item(113, 12).Value = 0  # L113
$ws.Cells.Item(113, 13).Value = -747.5  # M113
